$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for cells whose new values look like plain numbers,
# to preserve them as text (matching the original inlineStr string cells).
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '29.969.48'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '1.908.70'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '0.8125'
$ws.Range("E5").Value = '  +8.20%  '
$ws.Range("D6").Value = '241.37'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '0.3130'
$ws.Range("E8").Value = '  +3.24%  '
$ws.Range("D9").Value = '26.48'
$ws.Range("E9").Value = '  +4.53%  '
$ws.Range("E10").Value = '  +2.90%  '
$ws.Range("D11").Value = '0.08009'
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.906.11'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7447'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '5.185'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '92.53'
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").Value = '29.971.04'
$ws.Range("D17").Value = '14.01'
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").Value = '5.873'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '245.24'
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").Value = '0.000007787'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '2.152.04'
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").Value = '0.9993'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '6.942'
$ws.Range("E24").Value = '  +0.88%  '
$ws.Range("D25").Value = '0.1542'
$ws.Range("E25").Value = '  +21.50%  '
$ws.Range("D26").Value = '168.66'
$ws.Range("E26").Value = '  +1.72%  '
$ws.Range("D27").Value = '9.214'
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").Value = '2.070'
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("D30").Value = '1.361'
$ws.Range("E30").Value = '  -1.73%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").Value = '4.069'
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("D34").Value = '0.05521'
$ws.Range("E34").Value = '  +6.26%  '
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("D36").Value = '0.7302'
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").Value = '2.708'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '0.01922'
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").Value = '0.4411'
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").Value = '72.11'
$ws.Range("E41").Value = '  +1.42%  '
$ws.Range("D42").Value = '5.991'
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("D43").Value = '0.9991'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '0.8380'
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("D45").Value = '1.892'
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").Value = '100.91'
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("D47").Value = '7.574'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '9.717'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '982.96'
$ws.Range("E49").Value = '  +9.72%  '
$ws.Range("D50").Value = '2.057.29'
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '36.16'
$ws.Range("E51").Value = '  +0.86%  '
